# Add a new worksheet "Najvecd zadetkov po sezonah" at the end of the workbook
$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("Napredovanje_izpad")
$srcHeaderStyle = $srcSheet.Range("B1:D1")
$srcIndexStyle = $srcSheet.Range("A2")

$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$ws.Name = "Najvecd zadetkov po sezonah"

# Header row
$ws.Cells.Item(1, 2).Value = "Season"
$ws.Cells.Item(1, 3).Value = "Team"
$ws.Cells.Item(1, 4).Value = "Goals for"

# Data rows: index, season, team, goals for
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "1993-1994"
$ws.Cells.Item(2, 3).Value = "Barcelona"
$ws.Cells.Item(2, 4).Value = 91
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "1994-1995"
$ws.Cells.Item(3, 3).Value = "Real Madrid"
$ws.Cells.Item(3, 4).Value = 76
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "1995-1996"
$ws.Cells.Item(4, 3).Value = "Valencia"
$ws.Cells.Item(4, 4).Value = 77
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "1996-1997"
$ws.Cells.Item(5, 3).Value = "Barcelona"
$ws.Cells.Item(5, 4).Value = 102
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "1997-1998"
$ws.Cells.Item(6, 3).Value = "Atletico Madrid"
$ws.Cells.Item(6, 4).Value = 79
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "1998-1999"
$ws.Cells.Item(7, 3).Value = "Barcelona"
$ws.Cells.Item(7, 4).Value = 87
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "1999-2000"
$ws.Cells.Item(8, 3).Value = "Barcelona"
$ws.Cells.Item(8, 4).Value = 69
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "2000-2001"
$ws.Cells.Item(9, 3).Value = "Real Madrid"
$ws.Cells.Item(9, 4).Value = 81
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "2001-2002"
$ws.Cells.Item(10, 3).Value = "Real Madrid"
$ws.Cells.Item(10, 4).Value = 69
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "2002-2003"
$ws.Cells.Item(11, 3).Value = "Real Madrid"
$ws.Cells.Item(11, 4).Value = 86
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "2003-2004"
$ws.Cells.Item(12, 3).Value = "Real Madrid"
$ws.Cells.Item(12, 4).Value = 72
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "2004-2005"
$ws.Cells.Item(13, 3).Value = "Barcelona"
$ws.Cells.Item(13, 4).Value = 73
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "2005-2006"
$ws.Cells.Item(14, 3).Value = "Barcelona"
$ws.Cells.Item(14, 4).Value = 80
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "2006-2007"
$ws.Cells.Item(15, 3).Value = "Barcelona"
$ws.Cells.Item(15, 4).Value = 78
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "2007-2008"
$ws.Cells.Item(16, 3).Value = "Real Madrid"
$ws.Cells.Item(16, 4).Value = 84
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "2008-2009"
$ws.Cells.Item(17, 3).Value = "Barcelona"
$ws.Cells.Item(17, 4).Value = 105
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "2009-2010"
$ws.Cells.Item(18, 3).Value = "Real Madrid"
$ws.Cells.Item(18, 4).Value = 102
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "2010-2011"
$ws.Cells.Item(19, 3).Value = "Real Madrid"
$ws.Cells.Item(19, 4).Value = 102
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "2011-2012"
$ws.Cells.Item(20, 3).Value = "Real Madrid"
$ws.Cells.Item(20, 4).Value = 121
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "2012-2013"
$ws.Cells.Item(21, 3).Value = "Barcelona"
$ws.Cells.Item(21, 4).Value = 115
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "2013-2014"
$ws.Cells.Item(22, 3).Value = "Real Madrid"
$ws.Cells.Item(22, 4).Value = 104
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "2014-2015"
$ws.Cells.Item(23, 3).Value = "Real Madrid"
$ws.Cells.Item(23, 4).Value = 118
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "2015-2016"
$ws.Cells.Item(24, 3).Value = "Barcelona"
$ws.Cells.Item(24, 4).Value = 112
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "2016-2017"
$ws.Cells.Item(25, 3).Value = "Barcelona"
$ws.Cells.Item(25, 4).Value = 116
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "2017-2018"
$ws.Cells.Item(26, 3).Value = "Barcelona"
$ws.Cells.Item(26, 4).Value = 99
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "2018-2019"
$ws.Cells.Item(27, 3).Value = "Barcelona"
$ws.Cells.Item(27, 4).Value = 90
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "2019-2020"
$ws.Cells.Item(28, 3).Value = "Barcelona"
$ws.Cells.Item(28, 4).Value = 86
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "2020-2021"
$ws.Cells.Item(29, 3).Value = "Barcelona"
$ws.Cells.Item(29, 4).Value = 85
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "2021-2022"
$ws.Cells.Item(30, 3).Value = "Real Madrid"
$ws.Cells.Item(30, 4).Value = 80
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "2022-2023"
$ws.Cells.Item(31, 3).Value = "Real Madrid"
$ws.Cells.Item(31, 4).Value = 75
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "2023-2024"
$ws.Cells.Item(32, 3).Value = "Real Madrid"
$ws.Cells.Item(32, 4).Value = 6

# Apply the same formatting used elsewhere in the workbook: bold/bordered/
# centered header row, and bold/bordered/centered index column (style reused
# from the existing "Napredovanje_izpad" sheet so no new style is created).
$srcHeaderStyle.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)

$srcIndexStyle.Copy()
$ws.Range("A2:A32").PasteSpecial(-4122)

$excel.CutCopyMode = $false
$ws.Range("A1").Select()
